$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    # Force the value to be written as text, preserving formats like
    # "59.883.27", "4.41", or "0.0520" without Excel auto-converting
    # them into numbers (which would drop formatting / precision).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '59.883.27'
Set-TextValue $ws.Range('E2') '  +0.37%  '
Set-TextValue $ws.Range('D3') '2.420.27'
Set-TextValue $ws.Range('E3') '  +0.59%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '550.79'
Set-TextValue $ws.Range('E5') '  -0.43%  '
Set-TextValue $ws.Range('D6') '137.63'
Set-TextValue $ws.Range('E6') '  +0.60%  '
Set-TextValue $ws.Range('E7') '  -0.05%  '
Set-TextValue $ws.Range('D8') '0.588'
Set-TextValue $ws.Range('E8') '  +2.87%  '
Set-TextValue $ws.Range('E9') '  -1.82%  '
Set-TextValue $ws.Range('D10') '5.68'
Set-TextValue $ws.Range('E10') '  -2.65%  '
Set-TextValue $ws.Range('B11') 'Cardano'
Set-TextValue $ws.Range('C11') 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextValue $ws.Range('D11') '0.355'
Set-TextValue $ws.Range('E11') '  -2.58%  '
Set-TextValue $ws.Range('B12') 'TRON'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D12') '0.147'
Set-TextValue $ws.Range('E12') '  -1.66%  '
Set-TextValue $ws.Range('D13') '25.55'
Set-TextValue $ws.Range('E13') '  +3.49%  '
Set-TextValue $ws.Range('D14') '2.851.23'
Set-TextValue $ws.Range('E14') '  +0.58%  '
Set-TextValue $ws.Range('D15') '59.804.60'
Set-TextValue $ws.Range('E15') '  +0.37%  '
Set-TextValue $ws.Range('D16') '0.0000138'
Set-TextValue $ws.Range('E16') '  -1.16%  '
Set-TextValue $ws.Range('D17') '2.422.82'
Set-TextValue $ws.Range('E17') '  +0.60%  '
Set-TextValue $ws.Range('D18') '11.39'
Set-TextValue $ws.Range('E18') '  +1.09%  '
Set-TextValue $ws.Range('D19') '4.41'
Set-TextValue $ws.Range('E19') '  -0.20%  '
Set-TextValue $ws.Range('D20') '329.21'
Set-TextValue $ws.Range('E20') '  -2.22%  '
Set-TextValue $ws.Range('D21') '6.70'
Set-TextValue $ws.Range('E21') '  -4.59%  '
Set-TextValue $ws.Range('E22') '  +0.01%  '
Set-TextValue $ws.Range('D23') '66.62'
Set-TextValue $ws.Range('E23') '  +2.63%  '
Set-TextValue $ws.Range('E24') '  +1.21%  '
Set-TextValue $ws.Range('D25') '8.71'
Set-TextValue $ws.Range('E25') '  +3.08%  '
Set-TextValue $ws.Range('E26') '  -0.77%  '
Set-TextValue $ws.Range('E27') '  +0.14%  '
Set-TextValue $ws.Range('E28') '  +1.67%  '
Set-TextValue $ws.Range('E29') '  -1.56%  '
Set-TextValue $ws.Range('D30') '168.91'
Set-TextValue $ws.Range('E30') '  -1.35%  '
Set-TextValue $ws.Range('E31') '  -2.08%  '
Set-TextValue $ws.Range('D32') '18.72'
Set-TextValue $ws.Range('E32') '  -0.32%  '
Set-TextValue $ws.Range('E33') '  -0.34%  '
Set-TextValue $ws.Range('E35') '  +0.09%  '
Set-TextValue $ws.Range('E36') '  +0.11%  '
Set-TextValue $ws.Range('E37') '  -2.27%  '
Set-TextValue $ws.Range('E38') '  -2.40%  '
Set-TextValue $ws.Range('B39') 'PolygonEcosystemToken'
Set-TextValue $ws.Range('C39') 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D39') '0.411'
Set-TextValue $ws.Range('E39') '  -2.69%  '
Set-TextValue $ws.Range('B40') 'Bittensor'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range('D40') '313.80'
Set-TextValue $ws.Range('E40') '  +4.56%  '
Set-TextValue $ws.Range('B41') 'Filecoin'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range('D41') '3.67'
Set-TextValue $ws.Range('E41') '  -1.79%  '
Set-TextValue $ws.Range('B42') 'Aave'
Set-TextValue $ws.Range('C42') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D42') '138.97'
Set-TextValue $ws.Range('E42') '  -1.69%  '
Set-TextValue $ws.Range('B43') 'Stellar'
Set-TextValue $ws.Range('C43') 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D43') '0.0967'
Set-TextValue $ws.Range('E43') '  +0.55%  '
Set-TextValue $ws.Range('B44') 'Hedera'
Set-TextValue $ws.Range('C44') 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D44') '0.0520'
Set-TextValue $ws.Range('E44') '  -1.11%  '
Set-TextValue $ws.Range('B45') 'InjectiveProtocol'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws.Range('D45') '19.54'
Set-TextValue $ws.Range('E45') '  +1.99%  '
Set-TextValue $ws.Range('B46') 'Mantle'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue $ws.Range('D46') '0.579'
Set-TextValue $ws.Range('E46') '  +0.80%  '
Set-TextValue $ws.Range('B47') 'VeChain'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range('D47') '0.0225'
Set-TextValue $ws.Range('E47') '  -0.86%  '
Set-TextValue $ws.Range('B48') 'Polygon'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D48') '0.387'
Set-TextValue $ws.Range('E48') '  -2.99%  '
Set-TextValue $ws.Range('B49') 'EnergySwap'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D49') '17.65'
Set-TextValue $ws.Range('E49') '  -0.44%  '
Set-TextValue $ws.Range('B50') 'WhiteBITCoin'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D50') '11.05'
Set-TextValue $ws.Range('E50') '  +0.15%  '
Set-TextValue $ws.Range('B51') 'dogwifhat'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range('D51') '1.58'
Set-TextValue $ws.Range('E51') '  -0.93%  '
